# Updates cryptos list (prices in column D, 1h volume % in column E)
# Rows 15/16 also swap Coin name + Link (WrappedBTC <-> Litecoin reordered)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = '27.258.76'
$ws.Range("E2").Value = '  +0.01%  '

# Row 3 (Ethereum)
$ws.Range("D3").Value = '1.632.05'
$ws.Range("E3").Value = '  -1.12%  '

# Row 5 (BNB)
$ws.Range("D5").Value = '''216.08'
$ws.Range("E5").Value = '  -0.91%  '

# Row 6 (XRP)
$ws.Range("D6").Value = '''0.522'
$ws.Range("E6").Value = '  +1.58%  '

# Row 7 (USDC)
$ws.Range("E7").Value = '  -0.26%  '

# Row 8 (Cardano)
$ws.Range("D8").Value = '''0.256'
$ws.Range("E8").Value = '  -0.33%  '

# Row 9 (Dogecoin)
$ws.Range("E9").Value = '  -0.74%  '

# Row 10 (Solana)
$ws.Range("D10").Value = '''20.32'
$ws.Range("E10").Value = '  +1.21%  '

# Row 11 (TRON)
$ws.Range("E11").Value = '  -0.14%  '

# Row 12 (WrappedEther)
$ws.Range("D12").Value = '1.654.09'
$ws.Range("E12").Value = '  +0.04%  '

# Row 13 (Polkadot)
$ws.Range("D13").Value = '''4.13'
$ws.Range("E13").Value = '  -0.34%  '

# Row 14 (Polygon)
$ws.Range("D14").Value = '''0.544'
$ws.Range("E14").Value = '  +0.26%  '

# Row 15 (WrappedBTC)
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").Value = '''65.01'
$ws.Range("E15").Value = '  -3.88%  '

# Row 16 (Litecoin)
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '27.207.60'
$ws.Range("E16").Value = '  -0.14%  '

# Row 17 (ShibaInu)
$ws.Range("D17").Value = '0.0₃0734'
$ws.Range("E17").Value = '  -0.81%  '

# Row 18 (BitcoinCash)
$ws.Range("D18").Value = '''217.11'
$ws.Range("E18").Value = '  -1.45%  '

# Row 19 (Dai)
$ws.Range("E19").Value = '  -0.15%  '

# Row 20 (Chainlink)
$ws.Range("E20").Value = '  +1.28%  '

# Row 21 (Uniswap)
$ws.Range("E21").Value = '  -1.43%  '

# Row 22 (Toncoin)
$ws.Range("E22").Value = '  -4.06%  '

# Row 23 (Avalanche)
$ws.Range("D23").Value = '''9.12'
$ws.Range("E23").Value = '  -1.02%  '

# Row 24 (Monero)
$ws.Range("D24").Value = '''148.12'
$ws.Range("E24").Value = '  +0.88%  '

# Row 25 (BinanceUSD)
$ws.Range("E25").Value = '  -0.20%  '

# Row 26 (Cosmos)
$ws.Range("E26").Value = '  -3.54%  '

# Row 27 (Stellar)
$ws.Range("E27").Value = '  -0.57%  '

# Row 28 (EthereumClassic)
$ws.Range("D28").Value = '''15.59'
$ws.Range("E28").Value = '  -1.49%  '

# Row 29 (Hedera)
$ws.Range("D29").Value = '''0.0507'
$ws.Range("E29").Value = '  -0.41%  '

# Row 30 (PancakeSwap)
$ws.Range("E30").Value = '  -0.94%  '

# Row 31 (Filecoin)
$ws.Range("D31").Value = '''3.39'
$ws.Range("E31").Value = '  -0.11%  '

# Row 32 (InternetComputer(DFINITY))
$ws.Range("E32").Value = '  -1.03%  '

# Row 33 (Maker)
$ws.Range("D33").Value = '1.317.74'
$ws.Range("E33").Value = '  +4.56%  '

# Row 34 (LidoDAOToken)
$ws.Range("D34").Value = '''1.56'
$ws.Range("E34").Value = '  -1.75%  '

# Row 35 (HuobiToken)
$ws.Range("D35").Value = '''2.45'
$ws.Range("E35").Value = '  -0.27%  '

# Row 36 (VeChain)
$ws.Range("E36").Value = '  -1.72%  '

# Row 37 (ImmutableX)
$ws.Range("E37").Value = '  -1.25%  '

# Row 38 (ARBITRUM)
$ws.Range("D38").Value = '''0.846'
$ws.Range("E38").Value = '  +0.39%  '

# Row 39 (PaxDollar)
$ws.Range("E39").Value = '  -0.17%  '

# Row 40 (MXToken)
$ws.Range("E40").Value = '  +1.88%  '

# Row 41 (TrustWalletToken)
$ws.Range("D41").Value = '''0.800'
$ws.Range("E41").Value = '  -1.21%  '

# Row 42 (Aave)
$ws.Range("D42").Value = '''63.59'
$ws.Range("E42").Value = '  +2.41%  '

# Row 43 (RocketPoolETH)
$ws.Range("D43").Value = '1.768.97'
$ws.Range("E43").Value = '  -1.31%  '

# Row 44 (FraxShare)
$ws.Range("E44").Value = '  -4.50%  '

# Row 45 (Quant)
$ws.Range("D45").Value = '''90.76'
$ws.Range("E45").Value = '  -1.07%  '

# Row 46 (RenderToken)
$ws.Range("E46").Value = '  +0.26%  '

# Row 47 (BabyDogeCoin)
$ws.Range("E47").Value = '  -0.31%  '

# Row 48 (WEMIXToken)
$ws.Range("D48").Value = '''0.819'
$ws.Range("E48").Value = '  +21.38%  '

# Row 49 (Cronos)
$ws.Range("E49").Value = '  +0.42%  '

# Row 50 (EnergySwap)
$ws.Range("D50").Value = '''7.52'
$ws.Range("E50").Value = '  -2.14%  '

# Row 51 (Algorand)
$ws.Range("D51").Value = '''0.0959'
$ws.Range("E51").Value = '  -1.11%  '
